$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 290.58823
$ws.Range("J29").Value = 364
$ws.Range("L29").Value = 1092
$ws.Range("N29").Value = -1654
$ws.Range("H33").Value = 560.5454999999999
$ws.Range("I33").Value = 116.4
$ws.Range("K33").Value = 116.4
$ws.Range("M33").Value = 112.6
$ws.Range("H53").Value = 744.4706
$ws.Range("J53").Value = 1588.6666
$ws.Range("L53").Value = 1588.6666
$ws.Range("N53").Value = -2862.6666
$ws.Range("H58").Value = 200
$ws.Range("J58").Value = 200
$ws.Range("L58").Value = 600
$ws.Range("N58").Value = -900
$ws.Range("H103").Value = 799.6667
$ws.Range("I103").Value = 799.5
$ws.Range("K103").Value = 2398.5
$ws.Range("M103").Value = -1812.5
$ws.Range("H111").Value = 129
$ws.Range("I111").Value = 129
$ws.Range("K111").Value = 387
$ws.Range("M111").Value = 2680
$ws.Range("H116").Value = 9500
$ws.Range("I116").Value = 9500
$ws.Range("K116").Value = 9500
$ws.Range("M116").Value = -6058
$ws.Range("H132").Value = 3518.8
$ws.Range("I132").Value = 3518.8
$ws.Range("K132").Value = 10556.4
$ws.Range("M132").Value = -8026.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1383.3684
$ws.Range("I32").Value = 1349.1389
$ws.Range("K32").Value = 1349.1389
$ws.Range("M32").Value = -1062.1389
$ws.Range("H110").Value = 986.875
$ws.Range("I110").Value = 986.875
$ws.Range("K110").Value = 986.875
$ws.Range("M110").Value = 1058.125
$ws.Range("H132").Value = 850
$ws.Range("I132").Value = 850
$ws.Range("K132").Value = 2550
$ws.Range("M132").Value = -20

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 99995
$ws.Range("J2").Value = 99995
$ws.Range("L2").Value = 99995
$ws.Range("N2").Value = -100221
$ws.Range("H13").Value = 99995
$ws.Range("J13").Value = 99995
$ws.Range("L13").Value = 99995
$ws.Range("N13").Value = -100331
$ws.Range("H14").Value = 449.875
$ws.Range("J14").Value = 449.875
$ws.Range("L14").Value = 449.875
$ws.Range("N14").Value = -793.875
$ws.Range("H26").Value = 20235.5
$ws.Range("I26").Value = 20235.5
$ws.Range("K26").Value = 20235.5
$ws.Range("M26").Value = -19943.5
$ws.Range("H56").Value = 1500
$ws.Range("I56").Value = 1500
$ws.Range("K56").Value = 1500
$ws.Range("M56").Value = -761
$ws.Range("H96").Value = 17500
$ws.Range("I96").Value = 17500
$ws.Range("K96").Value = 17500
$ws.Range("M96").Value = -14754
$ws.Range("H97").Value = 27500
$ws.Range("I97").Value = 27500
$ws.Range("K97").Value = 27500
$ws.Range("M97").Value = -26509
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("N100").Value = 0
$ws.Range("H105").Value = 3083.3333
$ws.Range("J105").Value = 2166.6667
$ws.Range("L105").Value = 2166.6667
$ws.Range("N105").Value = -5660.6667
$ws.Range("H115").Value = 80000
$ws.Range("I115").Value = 80000
$ws.Range("K115").Value = 80000
$ws.Range("M115").Value = -78433
$ws.Range("H134").Value = 2570.5715
$ws.Range("I134").Value = 2570.5715
$ws.Range("K134").Value = 7711.7145
$ws.Range("M134").Value = -5176.7145
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("N135").Value = 0
$ws.Range("H140").Value = 99995
$ws.Range("J140").Value = 99995
$ws.Range("L140").Value = 99995
$ws.Range("N140").Value = -110355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 55.666668
$ws.Range("I5").Value = 33.5
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 33.5
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 78.5
$ws.Range("N5").Value = -324
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("H105").Value = 1253.2727
$ws.Range("I105").Value = 567
$ws.Range("K105").Value = 567
$ws.Range("M105").Value = 1180
$ws.Range("H107").Value = 839.4
$ws.Range("I107").Value = 732.6667
$ws.Range("K107").Value = 732.6667
$ws.Range("M107").Value = 1187.3333
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H132").Value = 3998.25
$ws.Range("I132").Value = 2997
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 8991
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -6461
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 794.73334
$ws.Range("I12").Value = 436.375
$ws.Range("J12").Value = 1204.2858
$ws.Range("K12").Value = 1309.125
$ws.Range("L12").Value = 3612.8574
$ws.Range("M12").Value = -1136.125
$ws.Range("N12").Value = -3958.8574
$ws.Range("H68").Value = 2645
$ws.Range("I68").Value = 2650
$ws.Range("J68").Value = 2641.6667
$ws.Range("K68").Value = 7950
$ws.Range("L68").Value = 7925.000100000001
$ws.Range("M68").Value = -7139
$ws.Range("N68").Value = -9547.000100000001
$ws.Range("H71").Value = 2645
$ws.Range("I71").Value = 2650
$ws.Range("J71").Value = 2641.6667
$ws.Range("K71").Value = 23850
$ws.Range("L71").Value = 23775.0003
$ws.Range("M71").Value = -19794
$ws.Range("N71").Value = -31887.0003
$ws.Range("H106").Value = 1000
$ws.Range("J106").Value = 1000
$ws.Range("L106").Value = 3000
$ws.Range("N106").Value = -4892
$ws.Range("H131").Value = 1291.6364
$ws.Range("I131").Value = 938.625
$ws.Range("K131").Value = 2815.875
$ws.Range("M131").Value = 2224.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1287.25
$ws.Range("I126").Value = 1049.6666
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3148.9998
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -678.9998000000001
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 3464.4285
$ws.Range("I132").Value = 3767
$ws.Range("J132").Value = 1649
$ws.Range("K132").Value = 11301
$ws.Range("L132").Value = 4947
$ws.Range("M132").Value = -8771
$ws.Range("N132").Value = -10007
$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3194.2778
$ws.Range("J46").Value = 3205.7354
$ws.Range("L46").Value = 3205.7354
$ws.Range("N46").Value = -3581.7354
$ws.Range("H92").Value = 50000
$ws.Range("I92").Value = 50000
$ws.Range("K92").Value = 50000
$ws.Range("M92").Value = -47504
$ws.Range("H103").Value = 29999
$ws.Range("J103").Value = 29999
$ws.Range("L103").Value = 29999
$ws.Range("N103").Value = -32343
$ws.Range("H109").Value = 27285
$ws.Range("J109").Value = 27285
$ws.Range("L109").Value = 27285
$ws.Range("N109").Value = -30059
$ws.Range("H132").Value = 3579.3333
$ws.Range("I132").Value = 3579.3333
$ws.Range("K132").Value = 10737.9999
$ws.Range("M132").Value = -8207.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 5001000.5
$ws.Range("J5").Value = 10000000
$ws.Range("L5").Value = 10000000
$ws.Range("N5").Value = -10000224
$ws.Range("H107").Value = 7000
$ws.Range("I107").Value = 7000
$ws.Range("K107").Value = 21000
$ws.Range("M107").Value = -19080
$ws.Range("H132").Value = 2422
$ws.Range("I132").Value = 2422
$ws.Range("K132").Value = 7266
$ws.Range("M132").Value = -4736
